$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BAB")

# ---------------------------------------------------------------------------
# Add "Bounds by piecewise McCormick n=10" block (rows 33-39)
# ---------------------------------------------------------------------------
$ws.Range("A33").Value = "Bounds by piecewise McCormick n=10"

$ws.Range("B2").Copy()
$ws.Range("B34").PasteSpecial(-4122)
$ws.Range("B34").Value = "Number of scenarios"

$ws.Range("B3").Copy()
$ws.Range("B35").PasteSpecial(-4122)
$ws.Range("B35").Value = 5
$ws.Range("C35").Value = -18380

$ws.Range("B4").Copy()
$ws.Range("B36").PasteSpecial(-4122)
$ws.Range("B36").Value = 10

$ws.Range("B5").Copy()
$ws.Range("B37").PasteSpecial(-4122)
$ws.Range("B37").Value = 20

$ws.Range("B6").Copy()
$ws.Range("B38").PasteSpecial(-4122)
$ws.Range("B38").Value = 40

$ws.Range("B7").Copy()
$ws.Range("B39").PasteSpecial(-4122)
$ws.Range("B39").Value = 120

# ---------------------------------------------------------------------------
# Add "Bounds by logarithm piecewise McCormick n=5" block (rows 42-48)
# ---------------------------------------------------------------------------
$ws.Range("A42").Value = "Bounds by logarithm piecewise McCormick n=5"

$ws.Range("B2").Copy()
$ws.Range("B43").PasteSpecial(-4122)
$ws.Range("B43").Value = "Number of scenarios"

$ws.Range("B3").Copy()
$ws.Range("B44").PasteSpecial(-4122)
$ws.Range("B44").Value = 5
$ws.Range("C44").Value = -18359

$ws.Range("B4").Copy()
$ws.Range("B45").PasteSpecial(-4122)
$ws.Range("B45").Value = 10

$ws.Range("B5").Copy()
$ws.Range("B46").PasteSpecial(-4122)
$ws.Range("B46").Value = 20

$ws.Range("B6").Copy()
$ws.Range("B47").PasteSpecial(-4122)
$ws.Range("B47").Value = 40

$ws.Range("B7").Copy()
$ws.Range("B48").PasteSpecial(-4122)
$ws.Range("B48").Value = 120

# ---------------------------------------------------------------------------
# Switch the active/selected sheet from "Presentation" to "BAB", updating the
# selection on BAB to D44 and clearing the old selection highlight.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("D44").Select()
